# "Generate Report for Handback" - localization-status.xlsx
#
# This script reflects that the two e2e source files (1178f3c0... and
# d5b3f2d5...) have now been handed back and are in sync with en-US:
#   - Status moves from "Ready for handoff" to
#     "Handed back: in sync with en-US" everywhere it is shown.
#   - Each locale sheet (zh-cn, de-de) gets the resolved "Latest Target
#     File" (a link back to the source .md), the generated "Latest
#     Handback File" (the roundtripped .xlf) and the "Latest Handback
#     DateTime" filled in for both rows.
#   - A couple of columns are widened so the newly-populated values are
#     readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# --- GitHub blob URLs behind the existing A2/A3 (source-file) hyperlinks ---
$url1178 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/07756e0ee203ef41db1134c32bbd483299de7bee/e2e/1178f3c0-419d-45a4-9e96-5d56ce45f7eb.md"
$urld5b3 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/07756e0ee203ef41db1134c32bbd483299de7bee/e2e/d5b3f2d5-c12b-463a-9de5-c64ba9e54407.md"

$name1178 = "1178f3c0-419d-45a4-9e96-5d56ce45f7eb.md"
$named5b3 = "d5b3f2d5-c12b-463a-9de5-c64ba9e54407.md"

# ---------------------------------------------------------------------
# 1. Status column: "Ready for handoff" -> "Handed back: in sync with en-US"
# ---------------------------------------------------------------------

$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. zh-cn sheet: fill in Latest Target File / Latest Handback File /
#    Latest Handback DateTime for both rows.
# ---------------------------------------------------------------------

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $url1178, "", "", $name1178) | Out-Null
$zhcn.Range("J2").Value = "1178f3c0-419d-45a4-9e96-5d56ce45f7eb.f6e8a4b90ece0b4932b77da11abb738fc0f83668.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-31 07:42:47"

$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $urld5b3, "", "", $named5b3) | Out-Null
$zhcn.Range("J3").Value = "d5b3f2d5-c12b-463a-9de5-c64ba9e54407.40218fad912361533d5353c09392e6c4e341b544.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-31 07:42:47"

# ---------------------------------------------------------------------
# 3. de-de sheet: same idea, different timestamp/file names.
# ---------------------------------------------------------------------

$dede.Hyperlinks.Add($dede.Range("I2"), $url1178, "", "", $name1178) | Out-Null
$dede.Range("J2").Value = "1178f3c0-419d-45a4-9e96-5d56ce45f7eb.f6e8a4b90ece0b4932b77da11abb738fc0f83668.de-de.xlf"
$dede.Range("K2").Value = "2016-08-31 07:43:06"

$dede.Hyperlinks.Add($dede.Range("I3"), $urld5b3, "", "", $named5b3) | Out-Null
$dede.Range("J3").Value = "d5b3f2d5-c12b-463a-9de5-c64ba9e54407.40218fad912361533d5353c09392e6c4e341b544.de-de.xlf"
$dede.Range("K3").Value = "2016-08-31 07:43:06"

# ---------------------------------------------------------------------
# 4. Widen columns that now hold the longer status / file-name / link text.
#    (Excel's ColumnWidth is in "characters"; it gets snapped to the
#    nearest 1/6th of a character internally, so we pick the input value
#    that lands on the desired stored width.)
# ---------------------------------------------------------------------

$wideStatusWidth = 29.16666667   # -> stored width 30 (was ~17.2)
$wideFileWidth   = 39.16666667   # -> stored width 40 (was ~18.7 / 21.7)

# (numeric column indices: E=5, F=6, C=3, I=9, J=10)
$overview.Columns.Item(5).ColumnWidth = $wideStatusWidth
$overview.Columns.Item(6).ColumnWidth = $wideStatusWidth

$zhcn.Columns.Item(3).ColumnWidth = $wideStatusWidth
$zhcn.Columns.Item(9).ColumnWidth = $wideFileWidth
$zhcn.Columns.Item(10).ColumnWidth = $wideFileWidth

$dede.Columns.Item(3).ColumnWidth = $wideStatusWidth
$dede.Columns.Item(9).ColumnWidth = $wideFileWidth
$dede.Columns.Item(10).ColumnWidth = $wideFileWidth

Write-Host "Applied handback report updates"
